# "filled in lorem ipsum, finished subsec names"
# Fill in the weekly journal entries for column J (rows 5-13) and tighten
# up the row heights to match the re-typeset timesheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New journal text, column J -------------------------------------------------
# Entered top-down for the first few rows, then bottom-up for the rest -
# mirrors the order the author actually typed them in (keeps shared-string
# table ordering identical to the source commit).
$ws.Range("J5").Value = "timesheet, meeting with kyle, group management"
$ws.Range("J6").Value = "intro page plan, connections polishing, complete plan done"
$ws.Range("J7").Value = "presentation"
$ws.Range("J13").Value = "finished intro page, all is left to link pages"
$ws.Range("J12").Value = "talked with the css leaders about progress, made an idea on how to connect the pages within the pages, talked to the css leaders about my plan"
$ws.Range("J11").Value = "travel day, just talked to kyle about some admin stuff"
$ws.Range("J10").Value = "intro page development, talking to the team members about deadlines, reminders, checking on progress"
$ws.Range("J9").Value = "checking on progress, started writing the introductions for topics"
$ws.Range("J8").Value = "travel day "

# --- Row height retouch -----------------------------------------------------
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 16
$ws.Rows.Item(5).RowHeight = 27
$ws.Rows.Item(6).RowHeight = 16
$ws.Rows.Item(7).RowHeight = 16
$ws.Rows.Item(8).RowHeight = 16
$ws.Rows.Item(9).RowHeight = 16
$ws.Rows.Item(10).RowHeight = 16
$ws.Rows.Item(11).RowHeight = 16
$ws.Rows.Item(12).RowHeight = 16
$ws.Rows.Item(13).RowHeight = 16
$ws.Rows.Item(14).RowHeight = 16
$ws.Rows.Item(15).RowHeight = 16
$ws.Rows.Item(16).RowHeight = 16
$ws.Rows.Item(17).RowHeight = 16
$ws.Rows.Item(18).RowHeight = 16
$ws.Rows.Item(19).RowHeight = 16

# --- View state: scroll frozen pane over to column D, select J8 ------------
$excel.ActiveWindow.FreezePanes = $false
$ws.Columns.Item(3).Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("J8").Select()
